# Auto-generated: updates cryptos list price/volume columns (and the
# Fetch.AI / Binance-PegBSC-USD row swap) to match the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.306.99"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.493.93"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'587.27"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'134.55"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("D7").Value = "3.493.95"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'7.21"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("D13").Value = "4.085.47"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "3.492.73"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "64.351.20"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "'25.17"
$ws.Range("E18").Value = "  -9.31%  "
$ws.Range("D19").Value = "'10.07"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "'13.68"
$ws.Range("E21").Value = "  -4.92%  "
$ws.Range("D22").Value = "'385.44"
$ws.Range("D23").Value = "'0.567"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "3.629.06"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'74.18"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.54"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("D32").Value = "'8.26"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").Value = "3.512.83"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("D37").Value = "'23.41"
$ws.Range("E37").Value = "  -2.25%  "
$ws.Range("D38").Value = "'5.27"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").Value = "'6.86"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").Value = "'161.31"
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("D42").Value = "'0.0781"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("D44").Value = "'25.47"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").Value = "'4.40"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "2.469.90"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'6.74"
$ws.Range("E51").Value = "  -2.17%  "
